# Resize/reposition the 8 picture placeholders on slide 1 of the 2x4
# picture-grid layout. The new sizes/positions replace the previous
# "magic number" EMU values (e.g. 3044952, 3813048, 6099048, 9144000)
# with values derived from the image's true aspect ratio.
#
# PowerPoint's COM object model expresses Shape.Left/Top/Width/Height in
# points (1 pt = 12700 EMU), so the target EMU values below are divided
# by 12700 before being assigned.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Picture 1 - top-left narrow picture (position unchanged, only resized)
$s.Shapes.Item(1).Left   = 0
$s.Shapes.Item(1).Top    = 0
$s.Shapes.Item(1).Width  = 1760570 / 12700
$s.Shapes.Item(1).Height = 3047923 / 12700

# Picture 2 - top row, 2nd column
$s.Shapes.Item(2).Left   = 3047923 / 12700
$s.Shapes.Item(2).Top    = 0
$s.Shapes.Item(2).Width  = 3047923 / 12700
$s.Shapes.Item(2).Height = 3047923 / 12700

# Picture 3 - top row, 3rd column
# (NB: 6095847/12700 in points form is given explicitly below -- rather
# than as a division expression -- to avoid a double-precision rounding
# hair that would otherwise land one EMU short of the target when
# PowerPoint converts the point value back to EMU on save.)
$s.Shapes.Item(3).Left   = 479.9879608559055
$s.Shapes.Item(3).Top    = 0
$s.Shapes.Item(3).Width  = 3047923 / 12700
$s.Shapes.Item(3).Height = 3047923 / 12700

# Picture 4 - top row, 4th column
$s.Shapes.Item(4).Left   = 9143771 / 12700
$s.Shapes.Item(4).Top    = 0
$s.Shapes.Item(4).Width  = 3047923 / 12700
$s.Shapes.Item(4).Height = 3047923 / 12700

# Picture 5 - bottom row, 1st column
$s.Shapes.Item(5).Left   = 0
$s.Shapes.Item(5).Top    = 3429000 / 12700
$s.Shapes.Item(5).Width  = 3047923 / 12700
$s.Shapes.Item(5).Height = 3047923 / 12700

# Picture 6 - bottom row, 2nd column
$s.Shapes.Item(6).Left   = 3047923 / 12700
$s.Shapes.Item(6).Top    = 3429000 / 12700
$s.Shapes.Item(6).Width  = 3047923 / 12700
$s.Shapes.Item(6).Height = 3047923 / 12700

# Picture 7 - bottom row, 3rd column (see note on Picture 3 above)
$s.Shapes.Item(7).Left   = 479.9879608559055
$s.Shapes.Item(7).Top    = 3429000 / 12700
$s.Shapes.Item(7).Width  = 3047923 / 12700
$s.Shapes.Item(7).Height = 3047923 / 12700

# Picture 8 - bottom row, 4th column
$s.Shapes.Item(8).Left   = 9143771 / 12700
$s.Shapes.Item(8).Top    = 3429000 / 12700
$s.Shapes.Item(8).Width  = 3047923 / 12700
$s.Shapes.Item(8).Height = 3047923 / 12700
